$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3772.0513
$ws.Range("I15").Value = 3772.0513
$ws.Range("K15").Value = 11316.1539
$ws.Range("M15").Value = -11147.1539
$ws.Range("H106").Value = 1712.6
$ws.Range("I106").Value = 1129.9231
$ws.Range("J106").Value = 5500
$ws.Range("K106").Value = 1129.9231
$ws.Range("L106").Value = 5500
$ws.Range("M106").Value = -498.9231
$ws.Range("N106").Value = -6762
$ws.Range("H137").Value = 2044791
$ws.Range("I137").Value = 6251831
$ws.Range("J137").Value = 5014
$ws.Range("K137").Value = 18755493
$ws.Range("L137").Value = 15042
$ws.Range("M137").Value = -18752943
$ws.Range("N137").Value = -20142
$ws.Range("H138").Value = 5885205.5
$ws.Range("I138").Value = 3727.3845
$ws.Range("J138").Value = 9526120
$ws.Range("K138").Value = 11182.1535
$ws.Range("L138").Value = 28578360
$ws.Range("M138").Value = -6042.1535
$ws.Range("N138").Value = -28588640

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 41752670
$ws.Range("I61").Value = 55613270
$ws.Range("J61").Value = 170883.33
$ws.Range("K61").Value = 55613270
$ws.Range("L61").Value = 170883.33
$ws.Range("M61").Value = -55613058
$ws.Range("N61").Value = -171307.33
$ws.Range("H136").Value = 41752670
$ws.Range("I136").Value = 55613270
$ws.Range("J136").Value = 170883.33
$ws.Range("K136").Value = 166839810
$ws.Range("L136").Value = 512649.99
$ws.Range("M136").Value = -166837260
$ws.Range("N136").Value = -517749.99

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 16537.375
$ws.Range("I86").Value = 20533.166
$ws.Range("J86").Value = 4550
$ws.Range("K86").Value = 20533.166
$ws.Range("L86").Value = 4550
$ws.Range("M86").Value = -19410.166
$ws.Range("N86").Value = -6796
$ws.Range("H89").Value = 16537.375
$ws.Range("I89").Value = 20533.166
$ws.Range("J89").Value = 4550
$ws.Range("K89").Value = 102665.83
$ws.Range("L89").Value = 22750
$ws.Range("M89").Value = -97049.83
$ws.Range("N89").Value = -33982
$ws.Range("H107").Value = 2840.1035
$ws.Range("I107").Value = 2393.5557
$ws.Range("K107").Value = 2393.5557
$ws.Range("M107").Value = -473.5556999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 73967.21000000001
$ws.Range("I132").Value = 2154.3
$ws.Range("J132").Value = 253499.5
$ws.Range("K132").Value = 6462.900000000001
$ws.Range("L132").Value = 760498.5
$ws.Range("M132").Value = -3932.900000000001
$ws.Range("N132").Value = -765558.5
$ws.Range("H134").Value = 20027.416
$ws.Range("I134").Value = 2083.6191
$ws.Range("J134").Value = 61896.277
$ws.Range("K134").Value = 6250.8573
$ws.Range("L134").Value = 185688.831
$ws.Range("M134").Value = -3715.8573
$ws.Range("N134").Value = -190758.831

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 545.56757
$ws.Range("I113").Value = 508.42856
$ws.Range("J113").Value = 554.23334
$ws.Range("K113").Value = 1525.28568
$ws.Range("L113").Value = 1662.70002
$ws.Range("M113").Value = 644.71432
$ws.Range("N113").Value = -6002.70002
$ws.Range("H131").Value = 930.0323
$ws.Range("I131").Value = 420
$ws.Range("J131").Value = 1028.1154
$ws.Range("K131").Value = 1260
$ws.Range("L131").Value = 3084.3462
$ws.Range("M131").Value = 3780
$ws.Range("N131").Value = -13164.3462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 107.8
$ws.Range("I2").Value = 100.333336
$ws.Range("J2").Value = 119
$ws.Range("K2").Value = 100.333336
$ws.Range("L2").Value = 119
$ws.Range("M2").Value = 12.666664
$ws.Range("N2").Value = -345
$ws.Range("H70").Value = 118766.664
$ws.Range("I70").Value = 256500
$ws.Range("J70").Value = 8580
$ws.Range("K70").Value = 256500
$ws.Range("L70").Value = 8580
$ws.Range("M70").Value = -256230
$ws.Range("N70").Value = -9120
$ws.Range("H73").Value = 118766.664
$ws.Range("I73").Value = 256500
$ws.Range("J73").Value = 8580
$ws.Range("K73").Value = 256500
$ws.Range("L73").Value = 8580
$ws.Range("M73").Value = -255564
$ws.Range("N73").Value = -10452
$ws.Range("H132").Value = 82290.75999999999
$ws.Range("I132").Value = 51913.8
$ws.Range("J132").Value = 203798.6
$ws.Range("K132").Value = 155741.4
$ws.Range("L132").Value = 611395.8
$ws.Range("M132").Value = -153211.4
$ws.Range("N132").Value = -616455.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 145.18182
$ws.Range("I55").Value = 145.18182
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 145.18182
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = 27.81818000000001
$ws.Range("H61").Value = 1889.409
$ws.Range("I61").Value = 1586.7778
$ws.Range("J61").Value = 3251.25
$ws.Range("K61").Value = 1586.7778
$ws.Range("L61").Value = 3251.25
$ws.Range("M61").Value = -1384.7778
$ws.Range("N61").Value = -3655.25
$ws.Range("H82").Value = 1741.5714
$ws.Range("I82").Value = 1345.5
$ws.Range("K82").Value = 1345.5
$ws.Range("M82").Value = -984.5
$ws.Range("H85").Value = 1741.5714
$ws.Range("I85").Value = 1345.5
$ws.Range("K85").Value = 1345.5
$ws.Range("M85").Value = -97.5
$ws.Range("H93").Value = 1563.909
$ws.Range("I93").Value = 1470.3
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 1470.3
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -222.3
$ws.Range("N93").Value = -4996
$ws.Range("H100").Value = 1324.4615
$ws.Range("I100").Value = 1035.8889
$ws.Range("K100").Value = 1035.8889
$ws.Range("M100").Value = -494.8888999999999
$ws.Range("H113").Value = 1889.409
$ws.Range("I113").Value = 1586.7778
$ws.Range("J113").Value = 3251.25
$ws.Range("K113").Value = 1586.7778
$ws.Range("L113").Value = 3251.25
$ws.Range("M113").Value = 583.2221999999999
$ws.Range("N113").Value = -7591.25
$ws.Range("H132").Value = 253499.5
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 337332.66
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 1011997.98
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -1017057.98
$ws.Range("H136").Value = 37569.07
$ws.Range("I136").Value = 22411.883
$ws.Range("J136").Value = 148000
$ws.Range("K136").Value = 67235.649
$ws.Range("L136").Value = 444000
$ws.Range("M136").Value = -64685.649
$ws.Range("N136").Value = -449100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 11094.857
$ws.Range("J45").Value = 10844.167
$ws.Range("L45").Value = 10844.167
$ws.Range("N45").Value = -11826.167
$ws.Range("H108").Value = 40625
$ws.Range("J108").Value = 40625
$ws.Range("L108").Value = 40625
$ws.Range("N108").Value = -48305
$ws.Range("H132").Value = 112823.72
$ws.Range("I132").Value = 91950.37
$ws.Range("J132").Value = 145624.72
$ws.Range("K132").Value = 275851.11
$ws.Range("L132").Value = 436874.16
$ws.Range("M132").Value = -273321.11
$ws.Range("N132").Value = -441934.16
$ws.Range("H136").Value = 46133.777
$ws.Range("I136").Value = 38260
$ws.Range("J136").Value = 57944.445
$ws.Range("K136").Value = 114780
$ws.Range("L136").Value = 173833.335
$ws.Range("M136").Value = -112230
$ws.Range("N136").Value = -178933.335
